$d = $word.ActiveDocument

$d.Content.Find.Execute("62×11=", $true, $false, $false, $false, $false, $true, 1, $false, "22×77=", 2) | Out-Null
$d.Content.Find.Execute("32×35=", $true, $false, $false, $false, $false, $true, 1, $false, "19×97=", 2) | Out-Null
$d.Content.Find.Execute("38×80=", $true, $false, $false, $false, $false, $true, 1, $false, "42×26=", 2) | Out-Null
$d.Content.Find.Execute("43×49=", $true, $false, $false, $false, $false, $true, 1, $false, "66×55=", 2) | Out-Null
$d.Content.Find.Execute("29×39=", $true, $false, $false, $false, $false, $true, 1, $false, "99×62=", 2) | Out-Null
$d.Content.Find.Execute("18×45=", $true, $false, $false, $false, $false, $true, 1, $false, "18×84=", 2) | Out-Null
$d.Content.Find.Execute("14×98=", $true, $false, $false, $false, $false, $true, 1, $false, "40×97=", 2) | Out-Null
$d.Content.Find.Execute("33×29=", $true, $false, $false, $false, $false, $true, 1, $false, "53×60=", 2) | Out-Null
$d.Content.Find.Execute("83×40=", $true, $false, $false, $false, $false, $true, 1, $false, "41×45=", 2) | Out-Null
$d.Content.Find.Execute("29×88=", $true, $false, $false, $false, $false, $true, 1, $false, "45×17=", 2) | Out-Null
$d.Content.Find.Execute("95×97=", $true, $false, $false, $false, $false, $true, 1, $false, "75×34=", 2) | Out-Null
$d.Content.Find.Execute("73×53=", $true, $false, $false, $false, $false, $true, 1, $false, "28×16=", 2) | Out-Null
$d.Content.Find.Execute("75×35=", $true, $false, $false, $false, $false, $true, 1, $false, "61×44=", 2) | Out-Null
$d.Content.Find.Execute("61×40=", $true, $false, $false, $false, $false, $true, 1, $false, "74×65=", 2) | Out-Null
$d.Content.Find.Execute("51×68=", $true, $false, $false, $false, $false, $true, 1, $false, "35×98=", 2) | Out-Null
$d.Content.Find.Execute("96×69=", $true, $false, $false, $false, $false, $true, 1, $false, "72×71=", 2) | Out-Null
$d.Content.Find.Execute("80×75=", $true, $false, $false, $false, $false, $true, 1, $false, "12×42=", 2) | Out-Null
$d.Content.Find.Execute("34×67=", $true, $false, $false, $false, $false, $true, 1, $false, "65×41=", 2) | Out-Null
$d.Content.Find.Execute("46×27=", $true, $false, $false, $false, $false, $true, 1, $false, "30×64=", 2) | Out-Null
$d.Content.Find.Execute("91×91=", $true, $false, $false, $false, $false, $true, 1, $false, "49×25=", 2) | Out-Null
$d.Content.Find.Execute("82×92=", $true, $false, $false, $false, $false, $true, 1, $false, "69×88=", 2) | Out-Null
$d.Content.Find.Execute("26×26=", $true, $false, $false, $false, $false, $true, 1, $false, "18×34=", 2) | Out-Null
$d.Content.Find.Execute("72×56=", $true, $false, $false, $false, $false, $true, 1, $false, "67×82=", 2) | Out-Null
$d.Content.Find.Execute("67×37=", $true, $false, $false, $false, $false, $true, 1, $false, "98×60=", 2) | Out-Null
$d.Content.Find.Execute("47×64=", $true, $false, $false, $false, $false, $true, 1, $false, "95×65=", 2) | Out-Null
